$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 330, shifting the existing rows 330-339 down to 332-341.
$ws.Rows("330:331").Insert()

# Fill in the new row 330 with the new weekly price record.
$ws.Range("A330").Value = 7
$ws.Range("B330").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C330").Value = "Ñuble"
$ws.Range("D330").Value = 45075
$ws.Range("E330").Value = 16
$ws.Range("F330").Value = "Fruta"
$ws.Range("G330").Value = 100104
$ws.Range("H330").Value = "Frutos de pepita"
$ws.Range("I330").Value = 100104005
$ws.Range("J330").Value = "Pera"
$ws.Range("K330").Value = "Forelle"
$ws.Range("L330").Value = "Especial"
$ws.Range("M330").Value = 60
$ws.Range("N330").Value = 12000
$ws.Range("O330").Value = 12000
$ws.Range("P330").Value = 12000
$ws.Range("Q330").Value = "$/bandeja 18 kilos granel"
$ws.Range("R330").Value = "Región de O'Higgins"
$ws.Range("S330").Value = 667
$ws.Range("T330").Value = 18

# Fill in the new row 331 with the new weekly price record.
$ws.Range("A331").Value = 7
$ws.Range("B331").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C331").Value = "Ñuble"
$ws.Range("D331").Value = 45075
$ws.Range("E331").Value = 16
$ws.Range("F331").Value = "Fruta"
$ws.Range("G331").Value = 100104
$ws.Range("H331").Value = "Frutos de pepita"
$ws.Range("I331").Value = 100104005
$ws.Range("J331").Value = "Pera"
$ws.Range("K331").Value = "Forelle"
$ws.Range("L331").Value = "Primera"
$ws.Range("M331").Value = 50
$ws.Range("N331").Value = 10000
$ws.Range("O331").Value = 10000
$ws.Range("P331").Value = 10000
$ws.Range("Q331").Value = "$/bandeja 18 kilos granel"
$ws.Range("R331").Value = "Región de O'Higgins"
$ws.Range("S331").Value = 556
$ws.Range("T331").Value = 18
